# Avance Migracion Sipro - subproducto + correcciones en otros componentes
$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Controllers sheet: progress updates
# -----------------------------------------------------------------------
$wsControllers = $wb.Worksheets.Item("Controllers")
$wsControllers.Range("D63").Value = 1
$wsControllers.Range("D86").Value = 0.05
$wsControllers.Range("D87").Value = 0.8
$wsControllers.Range("D88").Value = 0.05

# -----------------------------------------------------------------------
# Daos sheet: progress updates
# -----------------------------------------------------------------------
$wsDaos = $wb.Worksheets.Item("Daos")
$wsDaos.Range("C59").Value = 0.5
$wsDaos.Range("C95").Value = 0.05
$wsDaos.Range("C97").Value = 1
$wsDaos.Range("C98").Value = 0.7
$wsDaos.Range("C99").Value = 1

# -----------------------------------------------------------------------
# Vistas sheet: progress updates
# -----------------------------------------------------------------------
$wsVistas = $wb.Worksheets.Item("Vistas")
$wsVistas.Range("C39").Value = 1
$wsVistas.Range("C40").Value = 1
$wsVistas.Range("C41").Value = 1
$wsVistas.Range("C76").Value = 0.4
$wsVistas.Range("C77").Value = 1
$wsVistas.Range("C78").Value = 1

# -----------------------------------------------------------------------
# Avance sheet: manual "Programacion de Vista" progress value
# -----------------------------------------------------------------------
$wsAvance = $wb.Worksheets.Item("Avance")
$wsAvance.Range("D11").Value = 20.4

# -----------------------------------------------------------------------
# Restore the view / selection state recorded for each sheet.
# Selecting a range on a worksheet also activates that worksheet, so the
# sheet selected last becomes the active (tabSelected) sheet/tab.
# -----------------------------------------------------------------------
$wsAvance.Range("D12").Select()

$wsVistas.Range("C77:C78").Select()

$wsDaos.Range("C96").Select()

# Controllers must end up as the active sheet/tab.
$wsControllers.Activate()
$wsControllers.Range("D88").Select()
